$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "57.082.18"
$ws.Cells.Item(2, 5).Value = "  -1.29%  "

$ws.Cells.Item(3, 4).Value = "2.983.19"
$ws.Cells.Item(3, 5).Value = "  -2.29%  "

$ws.Cells.Item(4, 5).Value = "  +0.12%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "499.54"
$ws.Cells.Item(5, 5).Value = "  -4.97%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "137.51"
$ws.Cells.Item(6, 5).Value = "  -3.52%  "

$ws.Cells.Item(7, 5).Value = "  +0.09%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.429"
$ws.Cells.Item(8, 5).Value = "  -4.42%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "7.29"
$ws.Cells.Item(9, 5).Value = "  -4.67%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.107"
$ws.Cells.Item(10, 5).Value = "  -4.90%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.357"
$ws.Cells.Item(11, 5).Value = "  -4.31%  "

$ws.Cells.Item(12, 4).Value = "3.497.66"
$ws.Cells.Item(12, 5).Value = "  -2.25%  "

$ws.Cells.Item(13, 5).Value = "  -2.46%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "26.03"
$ws.Cells.Item(14, 5).Value = "  -3.75%  "

$ws.Cells.Item(15, 5).Value = "  -6.38%  "

$ws.Cells.Item(16, 4).Value = "57.145.38"
$ws.Cells.Item(16, 5).Value = "  -1.08%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "6.06"
$ws.Cells.Item(17, 5).Value = "  -3.63%  "

$ws.Cells.Item(18, 4).Value = "2.986.48"
$ws.Cells.Item(18, 5).Value = "  -1.91%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "12.60"
$ws.Cells.Item(19, 5).Value = "  -3.93%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "7.86"
$ws.Cells.Item(20, 5).Value = "  -3.84%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "320.82"
$ws.Cells.Item(21, 5).Value = "  -5.17%  "

$ws.Cells.Item(22, 5).Value = "  -0.05%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "5.72"
$ws.Cells.Item(23, 5).Value = "  +0.39%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "0.491"
$ws.Cells.Item(24, 5).Value = "  -2.28%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "62.99"
$ws.Cells.Item(25, 5).Value = "  -3.13%  "

$ws.Cells.Item(26, 5).Value = "  +0.26%  "

$ws.Cells.Item(27, 5).Value = "  -5.24%  "

$ws.Cells.Item(28, 4).Value = "0.0₃0891"
$ws.Cells.Item(28, 5).Value = "  -9.00%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "6.61"
$ws.Cells.Item(29, 5).Value = "  -4.96%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "7.11"
$ws.Cells.Item(30, 5).Value = "  -3.77%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.77"
$ws.Cells.Item(31, 5).Value = "  -4.50%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "1.16"
$ws.Cells.Item(32, 5).Value = "  -6.05%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "20.14"
$ws.Cells.Item(33, 5).Value = "  -4.80%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "154.62"
$ws.Cells.Item(34, 5).Value = "  -1.23%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "4.57"
$ws.Cells.Item(35, 5).Value = "  -3.85%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "5.77"
$ws.Cells.Item(36, 5).Value = "  -4.03%  "

$ws.Cells.Item(37, 5).Value = "  -7.02%  "

$ws.Cells.Item(38, 5).Value = "  -6.57%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.0664"
$ws.Cells.Item(39, 5).Value = "  -5.59%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "37.83"
$ws.Cells.Item(40, 5).Value = "  +0.24%  "

$ws.Cells.Item(41, 4).Value = "3.017.63"
$ws.Cells.Item(41, 5).Value = "  -2.37%  "

$ws.Cells.Item(42, 5).Value = "  +0.05%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "3.74"
$ws.Cells.Item(43, 5).Value = "  -3.82%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.644"
$ws.Cells.Item(44, 5).Value = "  -2.91%  "

$ws.Cells.Item(45, 4).Value = "2.191.15"
$ws.Cells.Item(45, 5).Value = "  -5.93%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "1.38"
$ws.Cells.Item(46, 5).Value = "  -6.53%  "

$ws.Cells.Item(47, 5).Value = "  -1.77%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.934"
$ws.Cells.Item(48, 5).Value = "  -9.64%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.0234"
$ws.Cells.Item(49, 5).Value = "  -5.00%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "19.18"
$ws.Cells.Item(50, 5).Value = "  -4.73%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.77"
$ws.Cells.Item(51, 5).Value = "  -11.78%  "
